$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, taken from the authoritative diff of the
# canonical OOXML (cryptos.xlsx GitHub Actions price/volume refresh).
$updates = @{
    'D2' = '73.361.34'
    'E2' = '  +1.72%  '
    'D3' = '4.057.93'
    'E3' = '  +0.80%  '
    'E4' = '  -0.19%  '
    'D5' = '566.64'
    'E5' = '  +5.36%  '
    'D6' = '151.57'
    'E6' = '  -0.18%  '
    'D7' = '4.051.36'
    'E7' = '  +0.83%  '
    'D8' = '0.698'
    'E8' = '  -0.32%  '
    'E9' = '  -0.06%  '
    'D10' = '0.769'
    'E10' = '  +2.60%  '
    'E11' = '  +0.57%  '
    'D12' = '54.38'
    'E12' = '  +14.22%  '
    'D13' = '0.0000329'
    'E13' = '  +1.25%  '
    'D14' = '11.12'
    'E14' = '  +4.00%  '
    'D15' = '4.703.71'
    'E15' = '  +0.66%  '
    'D16' = '4.052.29'
    'E16' = '  +0.36%  '
    'D17' = '14.55'
    'E17' = '  +3.38%  '
    'D18' = '20.86'
    'E18' = '  +1.52%  '
    'D19' = '1.22'
    'E19' = '  +2.53%  '
    'E20' = '  -0.39%  '
    'D21' = '73.105.11'
    'E21' = '  +1.49%  '
    'D22' = '448.57'
    'E22' = '  +4.54%  '
    'D23' = '98.10'
    'E23' = '  -0.25%  '
    'B24' = 'ImmutableX'
    'C24' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D24' = '3.57'
    'E24' = '  +1.46%  '
    'B25' = 'PancakeSwap'
    'C25' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D25' = '4.46'
    'E25' = '  +5.21%  '
    'D26' = '14.79'
    'E26' = '  +2.58%  '
    'D27' = '4.35'
    'E27' = '  +19.31%  '
    'D28' = '11.34'
    'E28' = '  +2.16%  '
    'D29' = '11.09'
    'E29' = '  +3.04%  '
    'E30' = '  +1.59%  '
    'D31' = '37.29'
    'E31' = '  +1.18%  '
    'D32' = '7.87'
    'E32' = '  +11.95%  '
    'D33' = '0.135'
    'E33' = '  +3.83%  '
    'E34' = '  +1.65%  '
    'D35' = '689.31'
    'E35' = '  +1.17%  '
    'D36' = '48.53'
    'E36' = '  +12.81%  '
    'D37' = '68.17'
    'E37' = '  +3.03%  '
    'D38' = '0.0₃0907'
    'E38' = '  +10.17%  '
    'E39' = '  +5.29%  '
    'D40' = '0.149'
    'E40' = '  -1.99%  '
    'E41' = '  -0.13%  '
    'B42' = 'Dai'
    'C42' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D42' = '0.999'
    'E42' = '  +0.04%  '
    'B43' = 'THORChain'
    'C43' = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
    'D43' = '11.15'
    'E43' = '  +15.22%  '
    'B44' = 'WEMIXToken'
    'C44' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D44' = '3.34'
    'E44' = '  -1.36%  '
    'B45' = 'VeChain'
    'C45' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D45' = '0.0496'
    'E45' = '  +2.05%  '
    'E46' = '  +0.05%  '
    'D47' = '0.153'
    'E47' = '  +1.68%  '
    'D48' = '2.69'
    'E48' = '  +3.19%  '
    'D49' = '3.57'
    'E49' = '  +7.65%  '
    'E50' = '  +4.48%  '
    'D51' = '3.31'
    'E51' = '  -1.41%  '
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $cell = $ws.Range($cellRef)

    if ($newValue -match '^-?\d+(\.\d+)?$') {
        # Values that look like plain numbers (e.g. "0.698", "54.38") must
        # stay text cells (matching the source inlineStr cells), otherwise
        # Excel's COM layer silently coerces them to numeric values and
        # mangles the exact display text (trailing zeros / scientific
        # notation). Forcing a text number format, assigning, then
        # clearing the format back off keeps the cell a plain text cell
        # with no leftover explicit style, same as the original.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.ClearFormats()
    } else {
        $cell.Value = $newValue
    }
}
